$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New matchup rows (Player_1, Points_1, Player_2, Points_2) for rows 291-313.
$data = @(
    @(6,7,7,13),
    @(5,15,4,5),
    @(3,8,1,12),
    @(4,3,6,17),
    @(6,13,5,7),
    @(4,6,6,14),
    @(4,8,3,12),
    @(5,17,4,3),
    @(4,16,6,4),
    @(5,8,2,12),
    @(4,5,3,15),
    @(7,6,5,14),
    @(2,12,4,8),
    @(5,12,4,8),
    @(4,13,5,7),
    @(6,6,8,14),
    @(6,19,5,1),
    @(3,5,1,15),
    @(2,2,5,18),
    @(3,15,3,5),
    @(4,13,4,7),
    @(9,13,3,7),
    @(5,14,4,6)
)

$startRow = 291
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Formula = "=B$row+D$row"
}

$lastRow = $startRow + $data.Count - 1

# Update the view to match where Excel landed after the new entries
# (new dimension is A1:E313, selection drops to the first empty row).
$excel.ActiveWindow.ScrollRow = 295
$ws.Range("A" + ($lastRow + 1)).Select()
